# Tripadvisor New Orleans shard 157 - workbook update
#
# 1) hotel_info: insert a new "State" column between "Hotel_Name" and "City",
#    populated with "Louisiana" for the existing data row.
# 2) Reorder the worksheet tabs so "review_info" comes before "hotel_info".

$wb = $excel.ActiveWorkbook

# --- 1. Insert the "State" column into hotel_info -------------------------
$wsHotel = $wb.Worksheets.Item("hotel_info")

# Column C currently holds "City" (A=STR, B=Hotel_Name, C=City, ...).
# Insert a fresh column there, shifting City (and everything after it) right,
# then fill in the new header + value.
$wsHotel.Columns("C:C").Insert()
$wsHotel.Cells.Item(1, 3).Value = "State"
$wsHotel.Cells.Item(2, 3).Value = "Louisiana"

# --- 2. Put review_info ahead of hotel_info in the tab order --------------
$wsReview = $wb.Worksheets.Item("review_info")
$wsReview.Move($wb.Worksheets.Item(1))
